# Apply the corrections described in the commit:
# "corrected species' name error centroura and ferox"
# The underlying fix restructures the DD (Data Deficient) species summary table:
#   - removes the erroneous "Albania" and "Germany" rows (their counts were
#     artifacts of the mis-assigned species records)
#   - inserts a new "NE" (Not Evaluated) status column after "Region"
#   - refreshes all remaining region figures (and the Total column) to the
#     corrected counts

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Germany" row (row 6) and "Albania" row (row 2).
# Delete Germany first so Albania's row index (2) is unaffected.
$ws.Rows.Item(6).Delete()
$ws.Rows.Item(2).Delete()

# Insert a new column B for the "NE" status, shifting DD..Total right one column.
$ws.Columns.Item(2).Insert()

# ----- Header row -----
$ws.Range("A1").Value = "Region"
$ws.Range("B1").Value = "NE"
$ws.Range("C1").Value = "DD"
$ws.Range("D1").Value = "LC"
$ws.Range("E1").Value = "NT"
$ws.Range("F1").Value = "VU"
$ws.Range("G1").Value = "EN"
$ws.Range("H1").Value = "CR"
$ws.Range("I1").Value = "Total"

# ----- Data rows -----
$data = @(
    @("Baltic Sea",        0, 0,  0,  0, 3, 2,  3,  8),
    @("Croatia",           0, 8,  3,  1, 1, 5,  1,  19),
    @("Europe",            0, 11, 17, 7, 6, 6,  6,  53),
    @("Ireland",           7, 0,  8,  12,4, 1,  2,  27),
    @("Italy",             0, 23, 7,  1, 0, 3,  4,  38),
    @("Mediterranean Sea", 0, 12, 8,  7, 7, 4,  14, 52),
    @("Norway",            0, 2,  0,  2, 0, 0,  1,  5),
    @("Sweden",            0, 0,  0,  0, 3, 1,  0,  4),
    @("World",             0, 2,  17, 9, 17,9,  8,  62)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $ws.Cells.Item($r, 9).Value = $row[8]
    $r = $r + 1
}
